# Replace "Autor" with "Author" in the header row (row 1) for the
# "Autor N" / "Autor N Institution" column pairs (D1:AO1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column pairs start at D (col 4) and step by 2: D/E, F/G, H/I, ... AN/AO
# Author numbers run from 1 to 19.
$startCol = 4
for ($i = 1; $i -le 19; $i++) {
    $authorCol = $startCol + ($i - 1) * 2
    $institutionCol = $authorCol + 1

    $ws.Cells.Item(1, $authorCol).Value = "Author $i"
    $ws.Cells.Item(1, $institutionCol).Value = "Author $i Institution"
}
